$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Covariate'
$ws.Range("B1").Value = 'Median (20%, 80%)'
$ws.Range("C1").Value = 'Mean (min, max)'
$ws.Range("D1").Value = 'Description'

$ws.Range("A2").Value = 'harvest'
$ws.Range("B2").Value = '0 (0–0.09)'
$ws.Range("C2").Value = '0.05 (0–0.92)'

$ws.Range("A3").Value = 'roads'
$ws.Range("B3").Value = '0 (0–0.01)'
$ws.Range("C3").Value = '0.01 (0–0.07)'

$ws.Range("A4").Value = 'veg_edges'
$ws.Range("B4").Value = '0.01 (0–0.02)'
$ws.Range("C4").Value = '0.01 (0–0.36)'

$ws.Range("A5").Value = 'wells'
$ws.Range("B5").Value = '0.01 (0–0.02)'
$ws.Range("C5").Value = '0.01 (0–0.3)'

$ws.Range("A6").Value = 'lc_grassland'
$ws.Range("B6").Value = '0.03 (0–0.07)'
$ws.Range("C6").Value = '0.04 (0–0.96)'

$ws.Range("A7").Value = 'lc_coniferous'
$ws.Range("B7").Value = '0.45 (0.22–0.79)'
$ws.Range("C7").Value = '0.48 (0–1)'

$ws.Range("A8").Value = 'lc_broadleaf'
$ws.Range("B8").Value = '0.1 (0.01–0.31)'
$ws.Range("C8").Value = '0.17 (0–1)'

$ws.Range("A9").Value = 'lc_mixed'
$ws.Range("B9").Value = '0.07 (0.01–0.17)'
$ws.Range("C9").Value = '0.11 (0–0.93)'

$ws.Range("A10").Value = 'lc_developed'
$ws.Range("B10").Value = '0.02 (0–0.06)'
$ws.Range("C10").Value = '0.04 (0–0.56)'

$ws.Range("A11").Value = 'lc_shrub'
$ws.Range("B11").Value = '0.06 (0.01–0.2)'
$ws.Range("C11").Value = '0.13 (0–1)'

$ws.Range("A12").Value = 'osm_industrial'
$ws.Range("B12").Value = '0 (0–0.02)'
$ws.Range("C12").Value = '0.02 (0–0.69)'

$ws.Range("A13").Value = 'natural_cohesion'
$ws.Range("B13").Value = '99.57 (99.08–99.82)'
$ws.Range("C13").Value = '99.32 (50.56–100)'

$ws.Range("A14").Value = 'anthropogenic_cohesion'
$ws.Range("B14").Value = '99.72 (99.11–99.88)'
$ws.Range("C14").Value = '98.89 (61.4–100)'

$ws.Range("A15").Value = 'natural_ed'
$ws.Range("B15").Value = '73.88 (44.89–131.19)'
$ws.Range("C15").Value = '98.49 (0–597.03)'

$ws.Range("A16").Value = 'anthropogenic_ed'
$ws.Range("B16").Value = '71.85 (43.4–128.41)'
$ws.Range("C16").Value = '96.57 (0–597.03)'

$ws.Range("A17").Value = 'natural_tca'
$ws.Range("B17").Value = '1351.06 (227.22–3562.92)'
$ws.Range("C17").Value = '1909.32 (0–7620.17)'

$ws.Range("A18").Value = 'anthropogenic_tca'
$ws.Range("B18").Value = '90.36 (10.23–390.43)'
$ws.Range("C18").Value = '229.66 (0–2662.77)'

$ws.Range("A19").Value = 'natural_cai_mn'
$ws.Range("B19").Value = '37.87 (23.02–58.18)'
$ws.Range("C19").Value = '40.95 (0–98.57)'

$ws.Range("A20").Value = 'anthro_cai_mn'
$ws.Range("B20").Value = '2.46 (0.85–7.42)'
$ws.Range("C20").Value = '6.72 (0–90.58)'

$ws.Range("A21").Value = 'forest_cai_mn'
$ws.Range("B21").Value = '39.37 (25.08–53.52)'
$ws.Range("C21").Value = '39.95 (0–97.8)'

$ws.Range("A22").Value = 'nonforest_cai_mn'
$ws.Range("B22").Value = '20.96 (8.12–40.1)'
$ws.Range("C22").Value = '24.6 (0–91.85)'

$ws.Range("A23").Value = 'nonveg_anthro_cai_mn'
$ws.Range("B23").Value = '4.01 (0.85–9.45)'
$ws.Range("C23").Value = '7.07 (0–87.3)'

$ws.Range("A24").Value = 'veg_anthro_cai_mn'
$ws.Range("B24").Value = '2.09 (0.69–4.36)'
$ws.Range("C24").Value = '3.49 (0–90.49)'

$ws.Range("A25").Value = 'forest_cohesion'
$ws.Range("B25").Value = '99.47 (98.93–99.76)'
$ws.Range("C25").Value = '99.2 (72.43–100)'

$ws.Range("A26").Value = 'nonforest_cohesion'
$ws.Range("B26").Value = '98.32 (96.83–99.13)'
$ws.Range("C26").Value = '97.72 (0–100)'

$ws.Range("A27").Value = 'nonveg_anthro_cohesion'
$ws.Range("B27").Value = '98.77 (97.43–99.43)'
$ws.Range("C27").Value = '98.03 (34.43–100)'

$ws.Range("A28").Value = 'veg_anthro_cohesion'
$ws.Range("B28").Value = '99.31 (98.24–99.7)'
$ws.Range("C28").Value = '98.39 (61.4–100)'

$ws.Range("A29").Value = 'forest_ed'
$ws.Range("B29").Value = '68.44 (45.78–108.18)'
$ws.Range("C29").Value = '91.01 (0–597.03)'

$ws.Range("A30").Value = 'nonforest_ed'
$ws.Range("B30").Value = '26.77 (8.86–50.99)'
$ws.Range("C30").Value = '33.94 (0–275.53)'

$ws.Range("A31").Value = 'nonveg_anthro_ed'
$ws.Range("B31").Value = '23.7 (6.29–47.93)'
$ws.Range("C31").Value = '31.41 (0–387.89)'

$ws.Range("A32").Value = 'veg_anthro_ed'
$ws.Range("B32").Value = '85.3 (49.61–154.26)'
$ws.Range("C32").Value = '109.24 (0–602.17)'

$ws.Range("A33").Value = 'forest_tca'
$ws.Range("B33").Value = '1031.99 (165.01–2888.4)'
$ws.Range("C33").Value = '1545.25 (0–7333.05)'

$ws.Range("A34").Value = 'nonforest_tca'
$ws.Range("B34").Value = '80.83 (8.08–393.89)'
$ws.Range("C34").Value = '274.52 (0–4339.95)'

$ws.Range("A35").Value = 'nonveg_anthro_tca'
$ws.Range("B35").Value = '6.65 (0.07–45.38)'
$ws.Range("C35").Value = '40.24 (0–2075.66)'

$ws.Range("A36").Value = 'veg_anthro_tca'
$ws.Range("B36").Value = '33.78 (1.78–208.87)'
$ws.Range("C36").Value = '138.77 (0–2137.73)'

$ws.Range("A37").Value = 'seismic'
$ws.Range("B37").Value = '0.01 (0–0.02)'
$ws.Range("C37").Value = '0.01 (0–0.1)'

$ws.Range("A38").Value = 'pipe_trans'
$ws.Range("B38").Value = '0.01 (0–0.03)'
$ws.Range("C38").Value = '0.02 (0–0.46)'
